$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col5a1"
$ws.Range("C2").Value = "Sdc3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7800876666666667
$ws.Range("H2").Value = 2.340263
$ws.Range("I2").Value = 0.004173077125706292
$ws.Range("J2").Value = 0.004173077125706292
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.12444933333333
$ws.Range("N2").Value = 78.373348
$ws.Range("O2").Value = 0.7238861157526749
$ws.Range("P2").Value = 0.7238861157526749
$ws.Range("Q2").Value = 20.37936072339156
$ws.Range("R2").Value = 183.414246510524
$ws.Range("S2").Value = 0.003020832591263864
$ws.Range("T2").Value = 0.003020832591263865

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col5a1"
$ws.Range("C3").Value = "Sdc3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.7800876666666667
$ws.Range("H3").Value = 2.340263
$ws.Range("I3").Value = 0.004173077125706292
$ws.Range("J3").Value = 0.004173077125706292
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.818542
$ws.Range("N3").Value = 11.455626
$ws.Range("O3").Value = 0.1058085282850919
$ws.Range("P3").Value = 0.1058085282850919
$ws.Range("Q3").Value = 2.978797518848667
$ws.Range("R3").Value = 26.809177669638
$ws.Range("S3").Value = 0.0004415471490911643
$ws.Range("T3").Value = 0.0004415471490911644

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col5a1"
$ws.Range("C4").Value = "Sdc3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.7800876666666667
$ws.Range("H4").Value = 2.340263
$ws.Range("I4").Value = 0.004173077125706292
$ws.Range("J4").Value = 0.004173077125706292
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.146179
$ws.Range("N4").Value = 18.438537
$ws.Range("O4").Value = 0.1703053559622332
$ws.Range("P4").Value = 0.1703053559622332
$ws.Range("Q4").Value = 4.794558435025667
$ws.Range("R4").Value = 43.151025915231
$ws.Range("S4").Value = 0.0007106973853512631
$ws.Range("T4").Value = 0.0007106973853512632

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col5a1"
$ws.Range("C5").Value = "Sdc3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 166.39918
$ws.Range("H5").Value = 499.19754
$ws.Range("I5").Value = 0.8901520194024567
$ws.Range("J5").Value = 0.8901520194024568
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.12444933333333
$ws.Range("N5").Value = 78.373348
$ws.Range("O5").Value = 0.7238861157526749
$ws.Range("P5").Value = 0.7238861157526749
$ws.Range("Q5").Value = 4347.086947018213
$ws.Range("R5").Value = 39123.78252316391
$ws.Range("S5").Value = 0.6443686877546441
$ws.Range("T5").Value = 0.6443686877546442

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col5a1"
$ws.Range("C6").Value = "Sdc3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 166.39918
$ws.Range("H6").Value = 499.19754
$ws.Range("I6").Value = 0.8901520194024567
$ws.Range("J6").Value = 0.8901520194024568
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.818542
$ws.Range("N6").Value = 11.455626
$ws.Range("O6").Value = 0.1058085282850919
$ws.Range("P6").Value = 0.1058085282850919
$ws.Range("Q6").Value = 635.4022575955601
$ws.Range("R6").Value = 5718.62031836004
$ws.Range("S6").Value = 0.09418567512297654
$ws.Range("T6").Value = 0.09418567512297656

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col5a1"
$ws.Range("C7").Value = "Sdc3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 166.39918
$ws.Range("H7").Value = 499.19754
$ws.Range("I7").Value = 0.8901520194024567
$ws.Range("J7").Value = 0.8901520194024568
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.146179
$ws.Range("N7").Value = 18.438537
$ws.Range("O7").Value = 0.1703053559622332
$ws.Range("P7").Value = 0.1703053559622332
$ws.Range("Q7").Value = 1022.71914573322
$ws.Range("R7").Value = 9204.47231159898
$ws.Range("S7").Value = 0.1515976565248361
$ws.Range("T7").Value = 0.1515976565248361

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col5a1"
$ws.Range("C8").Value = "Sdc3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 19.75417333333333
$ws.Range("H8").Value = 59.26251999999999
$ws.Range("I8").Value = 0.105674903471837
$ws.Range("J8").Value = 0.105674903471837
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.12444933333333
$ws.Range("N8").Value = 78.373348
$ws.Range("O8").Value = 0.7238861157526749
$ws.Range("P8").Value = 0.7238861157526749
$ws.Range("Q8").Value = 516.066900368551
$ws.Range("R8").Value = 4644.602103316959
$ws.Range("S8").Value = 0.07649659540676691
$ws.Range("T8").Value = 0.07649659540676694

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col5a1"
$ws.Range("C9").Value = "Sdc3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 19.75417333333333
$ws.Range("H9").Value = 59.26251999999999
$ws.Range("I9").Value = 0.105674903471837
$ws.Range("J9").Value = 0.105674903471837
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.818542
$ws.Range("N9").Value = 11.455626
$ws.Range("O9").Value = 0.1058085282850919
$ws.Range("P9").Value = 0.1058085282850919
$ws.Range("Q9").Value = 75.43214054861333
$ws.Range("R9").Value = 678.88926493752
$ws.Range("S9").Value = 0.01118130601302422
$ws.Range("T9").Value = 0.01118130601302422

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col5a1"
$ws.Range("C10").Value = "Sdc3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 19.75417333333333
$ws.Range("H10").Value = 59.26251999999999
$ws.Range("I10").Value = 0.105674903471837
$ws.Range("J10").Value = 0.105674903471837
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.146179
$ws.Range("N10").Value = 18.438537
$ws.Range("O10").Value = 0.1703053559622332
$ws.Range("P10").Value = 0.1703053559622332
$ws.Range("Q10").Value = 121.4126853036933
$ws.Range("R10").Value = 1092.71416773324
$ws.Range("S10").Value = 0.01799700205204583
$ws.Range("T10").Value = 0.01799700205204583
